$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 326.33334
$ws.Range("I2").Value = 382.66666
$ws.Range("K2").Value = 382.66666
$ws.Range("M2").Value = -269.66666
$ws.Range("H6").Value = 83390.664
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("H17").Value = 5264514
$ws.Range("J17").Value = 5264514
$ws.Range("L17").Value = 15793542
$ws.Range("N17").Value = -15793878
$ws.Range("H62").Value = 2717.4285
$ws.Range("I62").Value = 2730.6
$ws.Range("K62").Value = 2730.6
$ws.Range("M62").Value = -2106.6
$ws.Range("H64").Value = 6818.879
$ws.Range("J64").Value = 10633.2
$ws.Range("L64").Value = 10633.2
$ws.Range("N64").Value = -11129.2
$ws.Range("H65").Value = 2717.4285
$ws.Range("I65").Value = 2730.6
$ws.Range("K65").Value = 13653
$ws.Range("M65").Value = -10533
$ws.Range("H67").Value = 6818.879
$ws.Range("J67").Value = 10633.2
$ws.Range("L67").Value = 10633.2
$ws.Range("N67").Value = -12349.2
$ws.Range("H137").Value = 2302.5652
$ws.Range("I137").Value = 2113.1875
$ws.Range("K137").Value = 6339.5625
$ws.Range("M137").Value = -3789.5625
$ws.Range("H141").Value = 1498.6
$ws.Range("I141").Value = 1498.6
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 4495.799999999999
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 684.2000000000007
$ws.Range("N141").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 5222.1665
$ws.Range("I5").Value = 6068.8
$ws.Range("K5").Value = 6068.8
$ws.Range("M5").Value = -5956.8
$ws.Range("H32").Value = 8895.689
$ws.Range("I32").Value = 8895.689
$ws.Range("K32").Value = 8895.689
$ws.Range("M32").Value = -8608.689
$ws.Range("H45").Value = 103872.5
$ws.Range("I45").Value = 335575
$ws.Range("J45").Value = 4571.4287
$ws.Range("K45").Value = 335575
$ws.Range("L45").Value = 4571.4287
$ws.Range("M45").Value = -335198
$ws.Range("N45").Value = -5325.4287
$ws.Range("H57").Value = 15000
$ws.Range("I57").Value = 15000
$ws.Range("K57").Value = 15000
$ws.Range("M57").Value = -14516
$ws.Range("H135").Value = 141249.75
$ws.Range("J135").Value = 141249.75
$ws.Range("L135").Value = 141249.75
$ws.Range("N135").Value = -151389.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 5222.1665
$ws.Range("I4").Value = 6068.8
$ws.Range("K4").Value = 6068.8
$ws.Range("M4").Value = -5953.8
$ws.Range("H82").Value = 29796.666
$ws.Range("J82").Value = 42195
$ws.Range("L82").Value = 42195
$ws.Range("N82").Value = -42961
$ws.Range("H85").Value = 29796.666
$ws.Range("J85").Value = 42195
$ws.Range("L85").Value = 42195
$ws.Range("N85").Value = -44847
$ws.Range("H97").Value = 16680
$ws.Range("J97").Value = 20503
$ws.Range("L97").Value = 20503
$ws.Range("N97").Value = -22485
$ws.Range("H134").Value = 3375.3096
$ws.Range("I134").Value = 3335.6829
$ws.Range("K134").Value = 10007.0487
$ws.Range("M134").Value = -7472.048699999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5658.2
$ws.Range("I31").Value = 5296.4287
$ws.Range("J31").Value = 5853
$ws.Range("K31").Value = 5296.4287
$ws.Range("L31").Value = 5853
$ws.Range("M31").Value = -5001.4287
$ws.Range("N31").Value = -6443
$ws.Range("H34").Value = 5658.2
$ws.Range("I34").Value = 5296.4287
$ws.Range("J34").Value = 5853
$ws.Range("K34").Value = 5296.4287
$ws.Range("L34").Value = 5853
$ws.Range("M34").Value = -5094.4287
$ws.Range("N34").Value = -6257
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H58").Value = 6875.826
$ws.Range("I58").Value = 4795.769
$ws.Range("J58").Value = 9579.9
$ws.Range("K58").Value = 4795.769
$ws.Range("L58").Value = 9579.9
$ws.Range("M58").Value = -4592.769
$ws.Range("N58").Value = -9985.9
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H68").Value = 41000
$ws.Range("J68").Value = 41000
$ws.Range("L68").Value = 41000
$ws.Range("N68").Value = -42498
$ws.Range("H71").Value = 41000
$ws.Range("J71").Value = 41000
$ws.Range("L71").Value = 123000
$ws.Range("N71").Value = -130488
$ws.Range("H74").Value = 42341.6
$ws.Range("J74").Value = 42341.6
$ws.Range("L74").Value = 42341.6
$ws.Range("N74").Value = -44089.6
$ws.Range("H77").Value = 42341.6
$ws.Range("J77").Value = 42341.6
$ws.Range("L77").Value = 127024.8
$ws.Range("N77").Value = -135760.8
$ws.Range("H134").Value = 8723
$ws.Range("I134").Value = 6840
$ws.Range("J134").Value = 14999.667
$ws.Range("K134").Value = 20520
$ws.Range("L134").Value = 44999.001
$ws.Range("M134").Value = -17985
$ws.Range("N134").Value = -50069.001
$ws.Range("H136").Value = 6875.826
$ws.Range("I136").Value = 4795.769
$ws.Range("J136").Value = 9579.9
$ws.Range("K136").Value = 14387.307
$ws.Range("L136").Value = 28739.7
$ws.Range("M136").Value = -11837.307
$ws.Range("N136").Value = -33839.7
$ws.Range("H138").Value = 123473.81
$ws.Range("I138").Value = 84598.914
$ws.Range("J138").Value = 240098.5
$ws.Range("K138").Value = 84598.914
$ws.Range("L138").Value = 240098.5
$ws.Range("M138").Value = -79458.914
$ws.Range("N138").Value = -250378.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 599.8333
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H123").Value = 16599.5
$ws.Range("I123").Value = 4899.5
$ws.Range("J123").Value = 20499.5
$ws.Range("K123").Value = 14698.5
$ws.Range("L123").Value = 61498.5
$ws.Range("M123").Value = -12248.5
$ws.Range("N123").Value = -66398.5
$ws.Range("H131").Value = 19232954
$ws.Range("J131").Value = 2308.9583
$ws.Range("L131").Value = 6926.874899999999
$ws.Range("N131").Value = -17006.8749
$ws.Range("H138").Value = 1637.875
$ws.Range("I138").Value = 1637.875
$ws.Range("K138").Value = 4913.625
$ws.Range("M138").Value = 226.375

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 176465.39
$ws.Range("I113").Value = 224309.17
$ws.Range("K113").Value = 224309.17
$ws.Range("M113").Value = -222139.17
$ws.Range("H123").Value = 46449.7
$ws.Range("J123").Value = 46449.7
$ws.Range("L123").Value = 46449.7
$ws.Range("N123").Value = -51349.7
$ws.Range("H132").Value = 2487.8333
$ws.Range("I132").Value = 2046
$ws.Range("J132").Value = 9999
$ws.Range("K132").Value = 6138
$ws.Range("L132").Value = 29997
$ws.Range("M132").Value = -3608
$ws.Range("N132").Value = -35057

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 14556.833
$ws.Range("I132").Value = 19083.53
$ws.Range("J132").Value = 8637.308000000001
$ws.Range("K132").Value = 57250.59
$ws.Range("L132").Value = 25911.924
$ws.Range("M132").Value = -54720.59
$ws.Range("N132").Value = -30971.924
$ws.Range("H136").Value = 4606.4546
$ws.Range("J136").Value = 4285.375
$ws.Range("L136").Value = 12856.125
$ws.Range("N136").Value = -17956.125

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 17706.182
$ws.Range("I41").Value = 11000
$ws.Range("K41").Value = 11000
$ws.Range("M41").Value = -10610
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H86").Value = 83333.336
$ws.Range("I86").Value = 50000
$ws.Range("K86").Value = 50000
$ws.Range("M86").Value = -48877
$ws.Range("H89").Value = 83333.336
$ws.Range("I89").Value = 50000
$ws.Range("K89").Value = 250000
$ws.Range("M89").Value = -244384
